$d = $word.ActiveDocument

$replacements = @(
    @("376×7=2632", "624×8=4992"),
    @("587×7=4109", "757×5=3785"),
    @("388×4=1552", "731×8=5848"),
    @("398×3=1194", "398×8=3184"),
    @("988×9=8892", "939×3=2817"),
    @("145×3=435", "378×3=1134"),
    @("807×2=1614", "965×8=7720"),
    @("274×9=2466", "517×9=4653"),
    @("735×4=2940", "260×8=2080"),
    @("363×5=1815", "712×5=3560"),
    @("205×4=820", "202×6=1212"),
    @("409×4=1636", "437×4=1748"),
    @("446×9=4014", "693×9=6237"),
    @("629×9=5661", "115×8=920"),
    @("609×3=1827", "478×2=956"),
    @("901×8=7208", "564×6=3384"),
    @("230×3=690", "356×9=3204"),
    @("303×5=1515", "296×5=1480"),
    @("262×7=1834", "156×5=780"),
    @("138×5=690", "157×3=471"),
    @("153×6=918", "227×9=2043"),
    @("923×8=7384", "823×2=1646"),
    @("371×4=1484", "832×4=3328"),
    @("597×9=5373", "622×2=1244"),
    @("640×2=1280", "759×6=4554")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
